# Power supply voltage reduced
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple scalar input changes -------------------------------------------------
$ws.Range("B5").Value = 1010
$ws.Range("B7").Value = 10
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 7

# B15 used to compute 3.7/B11 ; now it is a fixed constant formula
$ws.Range("B15").Formula = "=0.5"

# A16 label changes from "V_ref_adc" to "Vcc" (supply voltage, was reference voltage)
$ws.Range("A16").Value = "Vcc"
$ws.Range("B16").Value = 5

# --- Expand the sweep table from 8 rows (19-26) to 12 rows (19-30) ---------------
# Insert 4 new (blank) rows so that existing rows 22,23,24,25,26 move down to
# 23,25,27,29,30 respectively, leaving fresh rows at 22,24,26,28 for new data points.
$ws.Rows("22:22").Insert()
$ws.Rows("24:24").Insert()
$ws.Rows("26:26").Insert()
$ws.Rows("28:28").Insert()

# --- Fix up column A --------------------------------------------------------------
# The newly inserted rows carry an empty column-A cell; they are all still part of
# the 100-Ohm sweep, so re-assert 100. The final two rows (previously A=100/A=50)
# change to A=20.
$ws.Range("A22").Value = 100
$ws.Range("A24").Value = 100
$ws.Range("A26").Value = 100
$ws.Range("A28").Value = 100
$ws.Range("A29").Value = 20
$ws.Range("A30").Value = 20

# --- Fill in column B (input current, uA) for every row of the table ------------
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 1
$ws.Range("B21").Value = 10
$ws.Range("B22").Value = 60
$ws.Range("B23").Value = 100
$ws.Range("B24").Value = 600
$ws.Range("B25").Value = 1000
$ws.Range("B26").Value = 6000
$ws.Range("B27").Value = 10000
$ws.Range("B28").Value = 20000
$ws.Range("B29").Value = 50000
$ws.Range("B30").Value = 100000

# --- Rebuild the formulas for every row of the table -----------------------------
# Column C: power dissipated uses the row's own resistance (A) instead of the
# fixed R_s_1 cell (B$1)
# Columns D-G: voltage-divider style gain computation now adds the supply
# voltage offset B$4
for ($r = 19; $r -le 30; $r++) {
    $ws.Range("C$r").Formula = "=B$r*A$r/1000/1000"
    $ws.Range("D$r").Formula = "=MIN(B`$16,C$r*B`$15*B`$5+B`$4)"
    $ws.Range("E$r").Formula = "=MIN(B`$16,C$r*B`$15*B`$6+B`$4)"
    $ws.Range("F$r").Formula = "=MIN(B`$16,C$r*B`$15*B`$7+B`$4)"
    $ws.Range("G$r").Formula = "=MIN(B`$16,C$r*B`$15*B`$8+B`$4)"
    $ws.Range("I$r").Formula = "=A$r*(B$r/1000000)^2*1000"
    $ws.Range("J$r").Formula = "=(B`$11-C$r-B`$12)*B$r/1000"
}
